$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value2 = "Datos actualizados a 6 de Abril de 2020 a las 09:22"

# --- Arabia Saudita (row 39): new cases / active cases updated ---
$ws.Cells.Item(39, 2).Value2 = 2463
$ws.Cells.Item(39, 3).Value2 = 61
$ws.Cells.Item(39, 5).Value2 = 1941

# --- Kuwait (row 78): active / recovered updated ---
$ws.Cells.Item(78, 4).Value2 = 103
$ws.Cells.Item(78, 5).Value2 = 452

# --- Letonia overtakes Bulgaria in the ranking (rows 80-81) ---
# The table is sorted by "Casos totales" (column B) descending. Letonia's
# refreshed total (542) now exceeds Bulgaria's (541), so Letonia's row
# moves above Bulgaria's: push the current row 80 (Bulgaria) down to row
# 81 unchanged, then write Letonia's new figures into row 80.
$name80 = $ws.Cells.Item(80, 1).Value2
$b80 = $ws.Cells.Item(80, 2).Value2
$c80 = $ws.Cells.Item(80, 3).Value2
$d80 = $ws.Cells.Item(80, 4).Value2
$e80 = $ws.Cells.Item(80, 5).Value2
$f80 = $ws.Cells.Item(80, 6).Value2
$g80 = $ws.Cells.Item(80, 7).Value2
$h80 = $ws.Cells.Item(80, 8).Value2

$ws.Cells.Item(81, 1).Value2 = $name80
$ws.Cells.Item(81, 2).Value2 = $b80
$ws.Cells.Item(81, 3).Value2 = $c80
$ws.Cells.Item(81, 4).Value2 = $d80
$ws.Cells.Item(81, 5).Value2 = $e80
$ws.Cells.Item(81, 6).Value2 = $f80
$ws.Cells.Item(81, 7).Value2 = $g80
$ws.Cells.Item(81, 8).Value2 = $h80

$ws.Cells.Item(80, 1).Value2 = "Letonia"
$ws.Cells.Item(80, 2).Value2 = 542
$ws.Cells.Item(80, 3).Value2 = 9
$ws.Cells.Item(80, 4).Value2 = 1
$ws.Cells.Item(80, 5).Value2 = 540
$ws.Cells.Item(80, 6).Value2 = 5
$ws.Cells.Item(80, 7).Value2 = 0
$ws.Cells.Item(80, 8).Value2 = 1

# --- Banglades jumps up past several countries (rows 123-128) ---
# Banglades' refreshed total (117) now exceeds Camboya (114), Paraguay
# (113), Trinidad yTobago (105), Ruanda (104) and Gibraltar (103), so its
# row moves up to position 123 (just after Guinea, 121). Every row from
# 123 to 127 shifts down one slot to make room; the former row 128
# (the stale Banglades entry, 88) falls off the bottom of this block.
for ($r = 127; $r -ge 123; $r--) {
    $name = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2

    $ws.Cells.Item($r + 1, 1).Value2 = $name
    $ws.Cells.Item($r + 1, 2).Value2 = $b
    $ws.Cells.Item($r + 1, 3).Value2 = $c
    $ws.Cells.Item($r + 1, 4).Value2 = $d
    $ws.Cells.Item($r + 1, 5).Value2 = $e
    $ws.Cells.Item($r + 1, 6).Value2 = $f
    $ws.Cells.Item($r + 1, 7).Value2 = $g
    $ws.Cells.Item($r + 1, 8).Value2 = $h
}

$ws.Cells.Item(123, 1).Value2 = "Banglades"
$ws.Cells.Item(123, 2).Value2 = 117
$ws.Cells.Item(123, 3).Value2 = 29
$ws.Cells.Item(123, 4).Value2 = 33
$ws.Cells.Item(123, 5).Value2 = 71
$ws.Cells.Item(123, 6).Value2 = 1
$ws.Cells.Item(123, 7).Value2 = 4
$ws.Cells.Item(123, 8).Value2 = 13
